# fix variable efficiency bug in data prep
#
# 1) "Definition" sheet: the node rows (7-18) were written in the wrong
#    order - realign each row's Object_Name with its correct position.
# 2) "Nodes" sheet: same misalignment - the per-node attributes
#    (balance_type / has_state / node_state_cap / frac_state_loss /
#    node_slack_penalty) were tied to the wrong Object_Name; fix the
#    names so the attributes sit on the right node.
# 3) "Object__to_from_node": the Electrolyzer's ordered_unit_flow_op
#    relationship pointed at unit__to_node / Hydrogen_Kasso; it should be
#    unit__from_node / Power_Kasso.
# 4) "Variable_Eff": the unit__node__node relationship referenced the
#    wrong node (Hydrogen_Kasso) and used wrong class (unit__to_node) -
#    fix to unit__from_node / Power_Kasso, and correct the operating
#    point efficiency values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Definition sheet - reorder Object_Name column (rows 7-18)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Definition")

$ws1.Cells.Item(7,1).Value  = "E-Methanol_storage_Kasso"
$ws1.Cells.Item(8,1).Value  = "E-Methanol_Kasso"
$ws1.Cells.Item(9,1).Value  = "Power_Kasso"
$ws1.Cells.Item(10,1).Value = "Raw_Methanol"
$ws1.Cells.Item(11,1).Value = "Waste_Heat"
$ws1.Cells.Item(12,1).Value = "Hydrogen_Kasso"
$ws1.Cells.Item(13,1).Value = "Power_Wholesale"
$ws1.Cells.Item(15,1).Value = "Carbon_Dioxide"
$ws1.Cells.Item(17,1).Value = "Vaporized_Carbon_Dioxide"
$ws1.Cells.Item(18,1).Value = "Water"

# ---------------------------------------------------------------------
# 2) Nodes sheet - realign names with their balance_type/has_state/etc.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Nodes")

# row 2 -> E-Methanol_storage_Kasso (balance_type_node, has_state, caps)
$ws2.Cells.Item(2,1).Value = "E-Methanol_storage_Kasso"
$ws2.Cells.Item(2,3).Value = "balance_type_node"
$ws2.Cells.Item(2,4).Value = "'true"
$ws2.Cells.Item(2,5).Value = 100000
$ws2.Cells.Item(2,6).Value = 0
$ws2.Cells.Item(2,7).Value = 100000

# row 3 -> E-Methanol_Kasso (balance_type/penalty unchanged)
$ws2.Cells.Item(3,1).Value = "E-Methanol_Kasso"

# row 4 -> Power_Kasso
$ws2.Cells.Item(4,1).Value = "Power_Kasso"
$ws2.Cells.Item(4,3).Value = "balance_type_node"
$ws2.Cells.Item(4,7).Value = 100000

# row 5 -> Raw_Methanol
$ws2.Cells.Item(5,1).Value = "Raw_Methanol"
$ws2.Cells.Item(5,7).Value = 100000

# row 6 -> Waste_Heat
$ws2.Cells.Item(6,1).Value = "Waste_Heat"
$ws2.Cells.Item(6,7).Value = ""

# row 7 -> Hydrogen_Kasso
$ws2.Cells.Item(7,1).Value = "Hydrogen_Kasso"
$ws2.Cells.Item(7,3).Value = "balance_type_node"
$ws2.Cells.Item(7,7).Value = 100000

# row 8 -> Power_Wholesale (balance_type_none, no state)
$ws2.Cells.Item(8,1).Value = "Power_Wholesale"
$ws2.Cells.Item(8,3).Value = "balance_type_none"
$ws2.Cells.Item(8,4).Value = ""
$ws2.Cells.Item(8,5).Value = ""
$ws2.Cells.Item(8,6).Value = ""
$ws2.Cells.Item(8,7).Value = ""

# row 9 -> District_Heating (unchanged)

# row 10 -> Carbon_Dioxide (balance_type_none)
$ws2.Cells.Item(10,1).Value = "Carbon_Dioxide"
$ws2.Cells.Item(10,3).Value = "balance_type_none"
$ws2.Cells.Item(10,7).Value = ""

# row 11 -> Hydrogen_storage_Kasso (unchanged)

# row 12 -> Vaporized_Carbon_Dioxide (balance_type/penalty unchanged)
$ws2.Cells.Item(12,1).Value = "Vaporized_Carbon_Dioxide"

# row 13 -> Water (balance_type_none)
$ws2.Cells.Item(13,1).Value = "Water"
$ws2.Cells.Item(13,3).Value = "balance_type_none"
$ws2.Cells.Item(13,7).Value = ""

# ---------------------------------------------------------------------
# 3) Object__to_from_node - fix Electrolyzer ordered_unit_flow_op row
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Object__to_from_node")
$ws5.Cells.Item(19,1).Value = "unit__from_node"
$ws5.Cells.Item(19,4).Value = "Power_Kasso"

# ---------------------------------------------------------------------
# 4) Variable_Eff - fix relationship class/node and efficiency values
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Variable_Eff")
$ws7.Cells.Item(1,3).Value = "unit__from_node"
$ws7.Cells.Item(3,3).Value = "Power_Kasso"

$ws7.Cells.Item(6,2).Value = 0.667
$ws7.Cells.Item(6,3).Value = 0.4
$ws7.Cells.Item(7,2).Value = 0.7
$ws7.Cells.Item(7,3).Value = 0.7
$ws7.Cells.Item(8,2).Value = 0.75
